$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all touched cells so numeric-looking strings
# (e.g. "0.999", "8.30", "67.233.76") are preserved verbatim as text,
# matching the original inlineStr (text) cell type used throughout columns B-E.
$textCells = @('D2', 'E2', 'D3', 'E3', 'D4', 'E4', 'D5', 'E5', 'D6', 'E6', 'D8', 'E8', 'E9', 'D10', 'E10', 'E11', 'D12', 'E12', 'D13', 'E13', 'D14', 'E14', 'E15', 'D16', 'E16', 'D17', 'E17', 'D18', 'E18', 'D19', 'E19', 'E20', 'D21', 'E21', 'E22', 'E23', 'D24', 'E24', 'D25', 'E25', 'D26', 'E26', 'D27', 'E27', 'E28', 'D29', 'E29', 'E30', 'E31', 'D32', 'E32', 'E33', 'D34', 'E34', 'D35', 'E35', 'D36', 'E36', 'B37', 'C37', 'D37', 'E37', 'B38', 'C38', 'D38', 'E38', 'D39', 'E39', 'B40', 'C40', 'D40', 'E40', 'B41', 'C41', 'D41', 'E41', 'E42', 'B43', 'C43', 'D43', 'E43', 'B44', 'C44', 'D44', 'E44', 'D45', 'E45', 'D46', 'E46', 'D47', 'E47', 'D48', 'E48', 'D50', 'E50', 'E51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '67.233.76'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').Value = '3.107.25'
$ws.Range('E3').Value = '  +0.59%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '574.44'
$ws.Range('E5').Value = '  -0.68%  '
$ws.Range('D6').Value = '178.62'
$ws.Range('E6').Value = '  +4.16%  '
$ws.Range('D8').Value = '3.105.19'
$ws.Range('E8').Value = '  +0.70%  '
$ws.Range('E9').Value = '  -0.79%  '
$ws.Range('D10').Value = '6.39'
$ws.Range('E10').Value = '  -0.68%  '
$ws.Range('E11').Value = '  +0.50%  '
$ws.Range('D12').Value = '0.468'
$ws.Range('E12').Value = '  -1.31%  '
$ws.Range('D13').Value = '0.0000242'
$ws.Range('E13').Value = '  -1.29%  '
$ws.Range('D14').Value = '36.21'
$ws.Range('E14').Value = '  -0.69%  '
$ws.Range('E15').Value = '  +0.14%  '
$ws.Range('D16').Value = '3.625.63'
$ws.Range('E16').Value = '  +0.65%  '
$ws.Range('D17').Value = '67.122.29'
$ws.Range('E17').Value = '  +0.16%  '
$ws.Range('D18').Value = '7.06'
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('D19').Value = '3.105.66'
$ws.Range('E19').Value = '  +0.19%  '
$ws.Range('E20').Value = '  +0.63%  '
$ws.Range('D21').Value = '490.62'
$ws.Range('E21').Value = '  +0.38%  '
$ws.Range('E22').Value = '  +0.27%  '
$ws.Range('E23').Value = '  -1.17%  '
$ws.Range('D24').Value = '83.63'
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').Value = '12.63'
$ws.Range('E25').Value = '  -2.63%  '
$ws.Range('D26').Value = '2.27'
$ws.Range('E26').Value = '  +0.54%  '
$ws.Range('D27').Value = '10.14'
$ws.Range('E27').Value = '  -2.58%  '
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('D29').Value = '7.96'
$ws.Range('E29').Value = '  +2.92%  '
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('E31').Value = '  -2.11%  '
$ws.Range('D32').Value = '28.22'
$ws.Range('E32').Value = '  +0.28%  '
$ws.Range('E33').Value = '  -0.31%  '
$ws.Range('D34').Value = '0.0₃0944'
$ws.Range('E34').Value = '  +0.21%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('D36').Value = '47.59'
$ws.Range('E36').Value = '  +2.38%  '
$ws.Range('B37').Value = 'Mantle'
$ws.Range('C37').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D37').Value = '0.948'
$ws.Range('E37').Value = '  -1.66%  '
$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D38').Value = '5.58'
$ws.Range('E38').Value = '  -2.80%  '
$ws.Range('D39').Value = '0.314'
$ws.Range('E39').Value = '  +2.85%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '2.02'
$ws.Range('E40').Value = '  +0.30%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').Value = '49.13'
$ws.Range('E41').Value = '  -1.24%  '
$ws.Range('E42').Value = '  +0.41%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').Value = '2.75'
$ws.Range('E43').Value = '  +7.97%  '
$ws.Range('B44').Value = 'Cosmos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D44').Value = '8.30'
$ws.Range('E44').Value = '  -1.37%  '
$ws.Range('D45').Value = '2.801.48'
$ws.Range('E45').Value = '  +0.49%  '
$ws.Range('D46').Value = '370.77'
$ws.Range('E46').Value = '  -2.57%  '
$ws.Range('D47').Value = '0.0346'
$ws.Range('E47').Value = '  -0.83%  '
$ws.Range('D48').Value = '135.97'
$ws.Range('E48').Value = '  +0.85%  '
$ws.Range('D50').Value = '25.56'
$ws.Range('E50').Value = '  +3.58%  '
$ws.Range('E51').Value = '  +4.15%  '
